$d = $word.ActiveDocument

# Replace the placeholder id text (merging the stray trailing-space run into
# the main run, and dropping the trailing space along with it).
$found = $d.Content.Find.Execute(
    "**ID__AFFARS_5333_topic_13__ID** ", $false, $false, $false, $false,
    $false, $true, 1, $false, "**ID__AFFARS_5333_291__ID**", 2)

# Give the first paragraph the same paragraph border spacing / indent that
# the rest of the body paragraphs already use.
$p1 = $d.Paragraphs.Item(1)
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

Write-Host "Updated id paragraph. Find replaced:" $found
